# Fill in the next batch of rows (28-42) on the grid worksheet:
#  - Column B: change "yes/no" placeholder to "yes"
#  - Column E: fill in the reference-to-test string
#  - Move the visible top-left cell / selection to reflect the newly
#    worked-on area of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 28; $row -le 42; $row++) {
    $ws.Range("B$row").Value = "yes"
    $ws.Range("E$row").Value = "HTBHasRoyalFlushDefinitions.java"
}

# Update the view: scroll down to show the newly edited rows and select E36:E42
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("E36:E42").Select()
